$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the underlying worksheet data
$ws.Range("D2").Value = 22.095992
$ws.Range("D3").Value = 94.358521

# Refresh the chart's cached series values so the embedded chart XML matches
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = $ws.Range("D2:D8")
